# drill down system tab clockin
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 31-40 (Cashable column C) as "T", matching the rest of the
# column which is already flagged.
for ($r = 31; $r -le 40; $r++) {
    $ws.Cells.Item($r, 3).Value = "T"
}

# Move the active selection to F38, matching the latest interaction point.
$ws.Range("F38").Select()
